$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.877458214759827
$ws.Range("B1").Value = 2.239890336990356
$ws.Range("C1").Value = 2.350868940353394
$ws.Range("D1").Value = 2.881145477294922
$ws.Range("E1").Value = 2.267235517501831
